# ADD results from server
# Update computed result values (row 2) on each year sheet with the
# latest values retrieved from the server.

$wb = $excel.ActiveWorkbook

# Sheet "2025"
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 1037.265132737054
$ws.Range("E2").Value = 28926.05393052954
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 16171.06685703679
$ws.Range("L2").Value = 48492.22142001599
$ws.Range("M2").Value = 10595.37713982
$ws.Range("N2").Value = 7015.544443014018
$ws.Range("O2").Value = 6978.613354318873

# Sheet "2030"
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 4157.588990853394
$ws.Range("E2").Value = 45991.90904307188
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 37079.12819938764
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 17449.04999683176
$ws.Range("N2").Value = 8950.626290977361
$ws.Range("O2").Value = 9689.183138434251

# Sheet "2035"
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 2754.31755456332
$ws.Range("B2").Value = 6368.910634126893
$ws.Range("E2").Value = 57457.45307013817
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 52465.73681402855
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 21912.87293902603
$ws.Range("N2").Value = 12929.21708841329
$ws.Range("O2").Value = 12821.53916790957

# Sheet "2040"
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 2754.31755456332
$ws.Range("B2").Value = 6368.910634126893
$ws.Range("E2").Value = 57457.45307013817
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 52465.73681402855
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 21912.87293902603
$ws.Range("N2").Value = 13044.31074660906
$ws.Range("O2").Value = 12821.53916790957

# Sheet "2045"
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 5713.151062849596
$ws.Range("B2").Value = 6368.910634126893
$ws.Range("E2").Value = 57457.45307013817
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 52465.73681402855
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 21912.87293902603
$ws.Range("N2").Value = 13482.79507252446
$ws.Range("O2").Value = 14873.28000436571

# Sheet "2050"
$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 5713.151062849596
$ws.Range("B2").Value = 6368.910634126893
$ws.Range("E2").Value = 57457.45307013817
$ws.Range("G2").Value = 8095.925712661834
$ws.Range("I2").Value = 52465.73681402855
$ws.Range("L2").Value = 54844.03303316472
$ws.Range("M2").Value = 21912.87293902603
$ws.Range("N2").Value = 13482.79507252446
$ws.Range("O2").Value = 14873.28000436571
